$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - columns keep same label set, just reindexed shared strings
$ws.Range("B1").Value = "fruits"
$ws.Range("C1").Value = "vegetables"
$ws.Range("D1").Value = "maize_cassava_beans"

# Establish shared-string insertion order to match target file:
# revenue(0), maize_cassava_beans(1), vegetables(2), fruits(3), share(4), irrigation_requirement(5)
$ws.Range("A4").Value = "share"

# Row 2 - irrigation_requirement (renamed from irrigation_volume, new values)
$ws.Range("A2").Value = "irrigation_requirement"
$ws.Range("B2").Value = 754.3
$ws.Range("C2").Value = 818.4
$ws.Range("D2").Value = 768.7

# Row 3 - revenue (new values, with number formatting on C3/D3)
$ws.Range("A3").Value = "revenue"
$ws.Range("B3").Value = 29394
$ws.Range("C3").Value = 105427.35042735044
$ws.Range("D3").Value = 5084.0157954935976
$ws.Range("C3:D3").NumberFormat = "0"

# Row 4 - share (new row) values
$ws.Range("B4").Value = 0.54008438818565396
$ws.Range("C4").Value = 0.21518987341772153
$ws.Range("D4").Value = 0.24472573839662448

# Column A width update (widened to fit "irrigation_requirement")
$ws.Columns.Item(1).ColumnWidth = 21

# Update selection to D3
$ws.Range("D3").Select()
